# Update automàtic: dades i banners [2026-02-14 22:50]
# Refreshes the per-station weather snapshot: extraction timestamps (col E),
# and the latest observation values that changed since the previous run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-14 22:48:26"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "83%"
$ws.Range("O2").Value = "-1.3 °C"
$ws.Range("E3").Value = "2026-02-14 22:48:29"
$ws.Range("L3").Value = "63.4 km/h - 251º 22:14 TU"
$ws.Range("N3").Value = "-8.7 °C 22:29 TU"
$ws.Range("O3").Value = "-5.4 °C"
$ws.Range("E4").Value = "2026-02-14 22:48:32"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "70%"
$ws.Range("J4").Value = "998.4 hPa"
$ws.Range("N4").Value = "5.0 °C 22:29 TU"
$ws.Range("O4").Value = "10.5 °C"
$ws.Range("E5").Value = "2026-02-14 22:48:34"
$ws.Range("N5").Value = "-8.2 °C 22:29 TU"
$ws.Range("O5").Value = "-5.4 °C"
$ws.Range("E6").Value = "2026-02-14 22:48:36"
$ws.Range("J6").Value = "998.4 hPa"
$ws.Range("E7").Value = "2026-02-14 22:48:39"
$ws.Range("J7").Value = "998.6 hPa"
$ws.Range("K7").Value = "13.7 MJ/m2"
$ws.Range("E8").Value = "2026-02-14 22:48:41"
$ws.Range("J8").Value = "998.4 hPa"
$ws.Range("E9").Value = "2026-02-14 22:48:44"
$ws.Range("O9").Value = "11.6 °C"
$ws.Range("E10").Value = "2026-02-14 22:48:47"
$ws.Range("E11").Value = "2026-02-14 22:48:48"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "58%"
$ws.Range("E12").Value = "2026-02-14 22:48:49"
$ws.Range("N12").Value = "9.2 °C 22:01 TU"
$ws.Range("E13").Value = "2026-02-14 22:48:50"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "66%"
$ws.Range("J13").Value = "1001.3 hPa"
$ws.Range("O13").Value = "3.8 °C"
$ws.Range("E14").Value = "2026-02-14 22:48:51"
$ws.Range("O14").Value = "13.3 °C"
$ws.Range("E15").Value = "2026-02-14 22:48:52"
$ws.Range("O15").Value = "11.1 °C"
$ws.Range("E16").Value = "2026-02-14 22:48:53"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "74%"
$ws.Range("E17").Value = "2026-02-14 22:48:54"
$ws.Range("L17").Value = "76.7 km/h - 358º 22:19 TU"
$ws.Range("N17").Value = "-1.3 °C 22:05 TU"
$ws.Range("O17").Value = "1.5 °C"
$ws.Range("E18").Value = "2026-02-14 22:48:55"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "72%"
$ws.Range("J18").Value = "998.6 hPa"
$ws.Range("L18").Value = "30.6 km/h - 69º 22:17 TU"
$ws.Range("E19").Value = "2026-02-14 22:48:56"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "73%"
$ws.Range("E20").Value = "2026-02-14 22:48:58"
$ws.Range("I20").Value = "5.1 mm"
$ws.Range("N20").Value = "-8.8 °C 22:29 TU"
$ws.Range("O20").Value = "-5.6 °C"
$ws.Range("E21").Value = "2026-02-14 22:48:59"
$ws.Range("J21").Value = "1001.0 hPa"
$ws.Range("E22").Value = "2026-02-14 22:49:01"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "84%"
$ws.Range("O22").Value = "-6.9 °C"
$ws.Range("E23").Value = "2026-02-14 22:49:04"
$ws.Range("I23").Value = "40.8 mm"
$ws.Range("N23").Value = "-9.2 °C 22:24 TU"
$ws.Range("O23").Value = "-6.3 °C"
$ws.Range("E24").Value = "2026-02-14 22:49:07"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "67%"
$ws.Range("J24").Value = "1002.7 hPa"
$ws.Range("K24").Value = "14.5 MJ/m2"
$ws.Range("E25").Value = "2026-02-14 22:49:09"
$ws.Range("I25").Value = "20.4 mm"
$ws.Range("E26").Value = "2026-02-14 22:49:11"
$ws.Range("E27").Value = "2026-02-14 22:49:14"
$ws.Range("E28").Value = "2026-02-14 22:49:16"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "63%"
$ws.Range("J28").Value = "998.3 hPa"
$ws.Range("L28").Value = "69.1 km/h - 274º 22:08 TU"
$ws.Range("E29").Value = "2026-02-14 22:49:19"
$ws.Range("E30").Value = "2026-02-14 22:49:21"
$ws.Range("J30").Value = "998.3 hPa"
$ws.Range("O30").Value = "11.4 °C"
$ws.Range("E31").Value = "2026-02-14 22:49:24"
$ws.Range("J31").Value = "997.5 hPa"
$ws.Range("N31").Value = "7.0 °C 22:11 TU"
$ws.Range("O31").Value = "9.1 °C"
$ws.Range("E32").Value = "2026-02-14 22:49:26"
$ws.Range("N32").Value = "1.6 °C 22:22 TU"
$ws.Range("O32").Value = "4.1 °C"
$ws.Range("E33").Value = "2026-02-14 22:49:29"
$ws.Range("J33").Value = "1000.7 hPa"
$ws.Range("E34").Value = "2026-02-14 22:49:31"
$ws.Range("N34").Value = "-5.3 °C 22:25 TU"
$ws.Range("E35").Value = "2026-02-14 22:49:34"
$ws.Range("J35").Value = "1005.1 hPa"
$ws.Range("N35").Value = "1.3 °C 22:06 TU"
$ws.Range("E36").Value = "2026-02-14 22:49:36"
$ws.Range("J36").Value = "999.1 hPa"
$ws.Range("N36").Value = "9.4 °C 22:03 TU"
$ws.Range("O36").Value = "11.8 °C"
$ws.Range("E37").Value = "2026-02-14 22:49:39"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "62%"
$ws.Range("J37").Value = "999.2 hPa"
$ws.Range("E38").Value = "2026-02-14 22:49:41"
$ws.Range("H38").NumberFormat = "@"
$ws.Range("H38").Value = "79%"
$ws.Range("N38").Value = "6.0 °C 22:14 TU"
$ws.Range("E39").Value = "2026-02-14 22:49:44"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "84%"
$ws.Range("I39").Value = "14.6 mm"
$ws.Range("N39").Value = "-8.8 °C 22:12 TU"
$ws.Range("O39").Value = "-6.0 °C"
$ws.Range("E40").Value = "2026-02-14 22:49:47"
$ws.Range("J40").Value = "1001.7 hPa"
$ws.Range("O40").Value = "7.1 °C"
$ws.Range("E41").Value = "2026-02-14 22:49:49"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "48%"
$ws.Range("J41").Value = "1000.4 hPa"
$ws.Range("N41").Value = "10.3 °C 22:17 TU"
$ws.Range("O41").Value = "13.2 °C"
$ws.Range("E42").Value = "2026-02-14 22:49:52"
$ws.Range("O42").Value = "11.6 °C"
$ws.Range("E43").Value = "2026-02-14 22:49:54"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "63%"
$ws.Range("E44").Value = "2026-02-14 22:49:56"
$ws.Range("I44").Value = "37.8 mm"
$ws.Range("N44").Value = "-8.4 °C 22:29 TU"
$ws.Range("O44").Value = "-5.6 °C"
$ws.Range("E45").Value = "2026-02-14 22:49:59"
$ws.Range("J45").Value = "1007.8 hPa"
$ws.Range("E46").Value = "2026-02-14 22:50:01"

# The "@" text format above pulls HUMITAT cells onto a new style; restore
# their original (bordered, general-format) cell style by pasting just the
# formatting back from the untouched neighboring cell in the same row.
$ws.Range("I2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("I11").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("I13").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("I16").Copy() | Out-Null
$ws.Range("H16").PasteSpecial(-4122) | Out-Null
$ws.Range("I18").Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4122) | Out-Null
$ws.Range("I19").Copy() | Out-Null
$ws.Range("H19").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("I24").Copy() | Out-Null
$ws.Range("H24").PasteSpecial(-4122) | Out-Null
$ws.Range("I28").Copy() | Out-Null
$ws.Range("H28").PasteSpecial(-4122) | Out-Null
$ws.Range("I37").Copy() | Out-Null
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("I38").Copy() | Out-Null
$ws.Range("H38").PasteSpecial(-4122) | Out-Null
$ws.Range("I39").Copy() | Out-Null
$ws.Range("H39").PasteSpecial(-4122) | Out-Null
$ws.Range("I41").Copy() | Out-Null
$ws.Range("H41").PasteSpecial(-4122) | Out-Null
$ws.Range("I43").Copy() | Out-Null
$ws.Range("H43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
